$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.63867
$ws.Cells.Item(2, 8).Value = 1.91601
$ws.Cells.Item(2, 9).Value = 0.02162993170452444
$ws.Cells.Item(2, 10).Value = 0.02162993170452444
$ws.Cells.Item(2, 13).Value = 10.045207
$ws.Cells.Item(2, 14).Value = 30.135621
$ws.Cells.Item(2, 15).Value = 0.9365108453707793
$ws.Cells.Item(2, 16).Value = 0.9365108453707794
$ws.Cells.Item(2, 17).Value = 6.415572354689999
$ws.Cells.Item(2, 18).Value = 57.74015119221
$ws.Cells.Item(2, 19).Value = 0.0202566656259164
$ws.Cells.Item(2, 20).Value = 0.02025666562591641
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.63867
$ws.Cells.Item(3, 8).Value = 1.91601
$ws.Cells.Item(3, 9).Value = 0.02162993170452444
$ws.Cells.Item(3, 10).Value = 0.02162993170452444
$ws.Cells.Item(3, 15).Value = 0.03971513502725754
$ws.Cells.Item(3, 16).Value = 0.03971513502725754
$ws.Cells.Item(3, 17).Value = 0.27206873642
$ws.Cells.Item(3, 18).Value = 2.44861862778
$ws.Cells.Item(3, 19).Value = 0.0008590356582755469
$ws.Cells.Item(3, 20).Value = 0.000859035658275547
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.63867
$ws.Cells.Item(4, 8).Value = 1.91601
$ws.Cells.Item(4, 9).Value = 0.02162993170452444
$ws.Cells.Item(4, 10).Value = 0.02162993170452444
$ws.Cells.Item(4, 13).Value = 0.255005
$ws.Cells.Item(4, 14).Value = 0.765015
$ws.Cells.Item(4, 15).Value = 0.02377401960196297
$ws.Cells.Item(4, 16).Value = 0.02377401960196297
$ws.Cells.Item(4, 17).Value = 0.16286404335
$ws.Cells.Item(4, 18).Value = 1.46577639015
$ws.Cells.Item(4, 19).Value = 0.0005142304203324842
$ws.Cells.Item(4, 20).Value = 0.0005142304203324843
$ws.Cells.Item(5, 9).Value = 0.9490095874171892
$ws.Cells.Item(5, 10).Value = 0.9490095874171893
$ws.Cells.Item(5, 13).Value = 10.045207
$ws.Cells.Item(5, 14).Value = 30.135621
$ws.Cells.Item(5, 15).Value = 0.9365108453707793
$ws.Cells.Item(5, 16).Value = 0.9365108453707794
$ws.Cells.Item(5, 17).Value = 281.4821496683659
$ws.Cells.Item(5, 18).Value = 2533.339347015294
$ws.Cells.Item(5, 19).Value = 0.8887577709770463
$ws.Cells.Item(5, 20).Value = 0.8887577709770466
$ws.Cells.Item(6, 9).Value = 0.9490095874171892
$ws.Cells.Item(6, 10).Value = 0.9490095874171893
$ws.Cells.Item(6, 15).Value = 0.03971513502725754
$ws.Cells.Item(6, 16).Value = 0.03971513502725754
$ws.Cells.Item(6, 19).Value = 0.03769004390643564
$ws.Cells.Item(6, 20).Value = 0.03769004390643564
$ws.Cells.Item(7, 9).Value = 0.9490095874171892
$ws.Cells.Item(7, 10).Value = 0.9490095874171893
$ws.Cells.Item(7, 13).Value = 0.255005
$ws.Cells.Item(7, 14).Value = 0.765015
$ws.Cells.Item(7, 15).Value = 0.02377401960196297
$ws.Cells.Item(7, 16).Value = 0.02377401960196297
$ws.Cells.Item(7, 17).Value = 7.145632297689999
$ws.Cells.Item(7, 18).Value = 64.31069067921
$ws.Cells.Item(7, 19).Value = 0.02256177253370704
$ws.Cells.Item(7, 20).Value = 0.02256177253370704
$ws.Cells.Item(8, 7).Value = 0.866931
$ws.Cells.Item(8, 8).Value = 2.600793
$ws.Cells.Item(8, 9).Value = 0.02936048087828625
$ws.Cells.Item(8, 10).Value = 0.02936048087828625
$ws.Cells.Item(8, 13).Value = 10.045207
$ws.Cells.Item(8, 14).Value = 30.135621
$ws.Cells.Item(8, 15).Value = 0.9365108453707793
$ws.Cells.Item(8, 16).Value = 0.9365108453707794
$ws.Cells.Item(8, 17).Value = 8.708501349717
$ws.Cells.Item(8, 18).Value = 78.376512147453
$ws.Cells.Item(8, 19).Value = 0.02749640876781646
$ws.Cells.Item(8, 20).Value = 0.02749640876781646
$ws.Cells.Item(9, 7).Value = 0.866931
$ws.Cells.Item(9, 8).Value = 2.600793
$ws.Cells.Item(9, 9).Value = 0.02936048087828625
$ws.Cells.Item(9, 10).Value = 0.02936048087828625
$ws.Cells.Item(9, 15).Value = 0.03971513502725754
$ws.Cells.Item(9, 16).Value = 0.03971513502725754
$ws.Cells.Item(9, 17).Value = 0.369306248506
$ws.Cells.Item(9, 18).Value = 3.323756236554
$ws.Cells.Item(9, 19).Value = 0.001166055462546352
$ws.Cells.Item(9, 20).Value = 0.001166055462546352
$ws.Cells.Item(10, 7).Value = 0.866931
$ws.Cells.Item(10, 8).Value = 2.600793
$ws.Cells.Item(10, 9).Value = 0.02936048087828625
$ws.Cells.Item(10, 10).Value = 0.02936048087828625
$ws.Cells.Item(10, 13).Value = 0.255005
$ws.Cells.Item(10, 14).Value = 0.765015
$ws.Cells.Item(10, 15).Value = 0.02377401960196297
$ws.Cells.Item(10, 16).Value = 0.02377401960196297
$ws.Cells.Item(10, 17).Value = 0.221071739655
$ws.Cells.Item(10, 18).Value = 1.989645656895
$ws.Cells.Item(10, 19).Value = 0.0006980166479234361
$ws.Cells.Item(10, 20).Value = 0.0006980166479234361
